$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style of the existing headers (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the Save column values for rows 2-14
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 1
    9 = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 0
}

foreach ($row in 2..14) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
